# Improve query handling; Add synonym mapping & fuzzy matching
#
# Adds a new Q&A pair ("What are sensors?") to the bottom of the table,
# matches the workbook's existing formatting conventions for the new
# row, re-applies the explicit row heights that the source file carries
# on every data row, and widens the Question/Answer columns (plus the
# orphaned column I width hint) to fit the content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new Question/Answer row -------------------------------------
$ws.Range("A52").Value = "What are sensors?"
$ws.Range("B52").Value = "A sensor is a device that detects changes in its environment and responds by generating a signal that can be measured and used for various purposes"

# Give the new row its own font styling (matches the theme-coloured Arial
# font used elsewhere once a cell is touched, while leaving the other
# rows' explicit header/body styles untouched).
$ws.Range("A52:B52").Font.ThemeColor = 1

# --- Restore per-row explicit height on every used row --------------------
$ws.Rows("1:67").RowHeight = 15.75

# --- Column sizing ---------------------------------------------------------
$ws.Columns("A").ColumnWidth = 118/3
$ws.Columns("B").ColumnWidth = 353/3
$ws.Columns("I").ColumnWidth = 85/3

# --- Selection / scroll position -------------------------------------------
[void]$ws.Range("B61").Select()
